$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure any D-column values that look like plain numbers are forced to remain
# text (matching the source data, which stores prices/labels as text strings).
$textCells = @("D5", "D7", "D14", "D15", "D18", "D20", "D22", "D25", "D28", "D31", "D32", "D33", "D35", "D39", "D42", "D43", "D44", "D48")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '25.895.34'
$ws.Range("E2").Value = '  +0.10%  '
$ws.Range("D3").Value = '1.589.44'
$ws.Range("E3").Value = '  -1.68%  '
$ws.Range("E4").Value = '  -0.11%  '
$ws.Range("D5").Value = '210.10'
$ws.Range("E5").Value = '  -1.14%  '
$ws.Range("E6").Value = '  -0.13%  '
$ws.Range("D7").Value = '0.485'
$ws.Range("E7").Value = '  -2.98%  '
$ws.Range("E8").Value = '  -0.26%  '
$ws.Range("E9").Value = '  +0.40%  '
$ws.Range("E10").Value = '  -0.41%  '
$ws.Range("E11").Value = '  -0.26%  '
$ws.Range("D12").Value = '1.811.11'
$ws.Range("E12").Value = '  -1.62%  '
$ws.Range("D13").Value = '1.580.88'
$ws.Range("E13").Value = '  -2.11%  '
$ws.Range("D14").Value = '4.04'
$ws.Range("E14").Value = '  -2.35%  '
$ws.Range("D15").Value = '0.512'
$ws.Range("E15").Value = '  -2.08%  '
$ws.Range("D16").Value = '25.902.61'
$ws.Range("E16").Value = '  +0.12%  '
$ws.Range("B17").Value = 'ShibaInu'
$ws.Range("C17").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D17").Value = '0.0₃0724'
$ws.Range("E17").Value = '  -1.31%  '
$ws.Range("B18").Value = 'Litecoin'
$ws.Range("C18").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D18").Value = '60.24'
$ws.Range("E18").Value = '  -1.84%  '
$ws.Range("E19").Value = '  -0.16%  '
$ws.Range("D20").Value = '193.52'
$ws.Range("E20").Value = '  +1.31%  '
$ws.Range("E21").Value = '  -0.60%  '
$ws.Range("D22").Value = '9.42'
$ws.Range("E22").Value = '  -0.77%  '
$ws.Range("E23").Value = '  -1.08%  '
$ws.Range("E24").Value = '  -1.38%  '
$ws.Range("D25").Value = '141.57'
$ws.Range("E25").Value = '  -1.44%  '
$ws.Range("E26").Value = '  -0.18%  '
$ws.Range("E27").Value = '  -0.28%  '
$ws.Range("D28").Value = '15.16'
$ws.Range("E28").Value = '  -0.32%  '
$ws.Range("E29").Value = '  -2.42%  '
$ws.Range("E30").Value = '  -5.30%  '
$ws.Range("D31").Value = '0.0474'
$ws.Range("E31").Value = '  -0.43%  '
$ws.Range("D32").Value = '3.13'
$ws.Range("E32").Value = '  +0.22%  '
$ws.Range("D33").Value = '3.05'
$ws.Range("E33").Value = '  -1.50%  '
$ws.Range("E34").Value = '  +1.25%  '
$ws.Range("D35").Value = '2.36'
$ws.Range("E35").Value = '  -2.05%  '
$ws.Range("D36").Value = '1.109.41'
$ws.Range("E36").Value = '  -1.48%  '
$ws.Range("E37").Value = '  -0.21%  '
$ws.Range("E38").Value = '  -1.41%  '
$ws.Range("D39").Value = '0.507'
$ws.Range("E39").Value = '  -0.45%  '
$ws.Range("E40").Value = '  -1.42%  '
$ws.Range("E41").Value = '  -6.05%  '
$ws.Range("D42").Value = '0.815'
$ws.Range("E42").Value = '  +8.86%  '
$ws.Range("D43").Value = '5.18'
$ws.Range("E43").Value = '  +2.75%  '
$ws.Range("D44").Value = '93.73'
$ws.Range("E44").Value = '  -4.51%  '
$ws.Range("D45").Value = '1.724.07'
$ws.Range("E45").Value = '  -1.61%  '
$ws.Range("D46").Value = '0.0₆0111'
$ws.Range("E46").Value = '  -1.78%  '
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("D48").Value = '53.57'
$ws.Range("E48").Value = '  -0.66%  '
$ws.Range("E49").Value = '  -1.53%  '
$ws.Range("E50").Value = '  -0.52%  '
$ws.Range("E51").Value = '  -0.18%  '

# Reset the style of the forced-text cells back to Normal so no stray
# number-format style is left applied (keeps formatting identical to source).
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
